# Append a new data row (row 4) to the "Artfynd" sheet, mirroring the
# structure of the existing rows (2 and 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$row = 4

# --- Plain numeric cells -------------------------------------------------
$ws.Cells.Item($row, 1).Value  = 112159375        # A  Id
$ws.Cells.Item($row, 2).Value  = 90689             # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value  = 5966              # E  TaxonId
$ws.Cells.Item($row, 17).Value = 560630            # Q  Ost
$ws.Cells.Item($row, 18).Value = 6607020           # R  Nord
$ws.Cells.Item($row, 19).Value = 10                # S  Noggrannhet

# --- Plain text cells ------------------------------------------------------
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"                         # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value  = "NT"                                  # D  Rödlistade
$ws.Cells.Item($row, 6).Value  = "Motaggsvamp"                         # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Sarcodon squamosus"                  # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Schaeff.) Quél."                    # H  Auktor
$ws.Cells.Item($row, 10).Value = "fruktkroppar"                        # J  Enhet
$ws.Cells.Item($row, 11).Value = "teleomorf"                           # K  Ålder-Stadium
$ws.Cells.Item($row, 16).Value = "Tränstagården, VNV, Vstm"            # P  Lokalnamn
$ws.Cells.Item($row, 20).Value = "Västmanland"                         # T  Län
$ws.Cells.Item($row, 21).Value = "Köping"                              # U  Kommun
$ws.Cells.Item($row, 22).Value = "Västmanland"                         # V  Provins
$ws.Cells.Item($row, 23).Value = "Munktorp"                            # W  Församling
$ws.Cells.Item($row, 35).Value = "äldre, gles tallskog; vid elljusspår"# AI Biotop-beskrivning
$ws.Cells.Item($row, 49).Value = "Henrik Berg"                         # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Henrik Berg"                         # AX Observatörer

# --- Text cells that look numeric/date: force text via leading apostrophe,
#     the same way Excel itself records a quote-prefixed entry, so they are
#     stored as text (matching the source data) instead of being
#     reinterpreted as a number or a date serial. ------------------------
$ws.Cells.Item($row, 9).Value  = "'1"                 # I  Antal
$ws.Cells.Item($row, 25).Value = "'2023-09-17"         # Y  Startdatum
$ws.Cells.Item($row, 27).Value = "'2023-09-17"         # AA Slutdatum

# --- Boolean cells ---------------------------------------------------------
$ws.Cells.Item($row, 30).Value = $false                # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false                # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false                # AG Ospontan

# --- Present-but-empty text cells (mirrors the empty inline strings used
#     for unfilled columns on the other rows). ---------------------------
$ws.Cells.Item($row, 14).Value = "'"                   # N  Metod
$ws.Cells.Item($row, 32).Value = "'"                   # AF Bestämningsmetod
$ws.Cells.Item($row, 46).Value = "'"                   # AT Bestämningsår
$ws.Cells.Item($row, 51).Value = "'"                   # AY Projektnamn
